$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# A handful of the new prices read as plain numbers (e.g. "6.56"), but this
# column stores its values as literal text (t="inlineStr" in the original),
# matching the source feed formatting (note values like "60.521.29" with
# multiple separators). A direct `.Value = "6.56"` assignment would let Excel
# auto-convert the text to a real number, so instead build the text in a
# scratch cell via a formula ("=""6.56""" always evaluates to the text string
# 6.56), copy it, and paste-special *values only* into the target cell. A
# pasted string value is stored as text without ever touching the target
# cell's number format/style, so it lands exactly like the original content.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""60.521.29"""
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$scratch.Formula = "=""2.585.22"""
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$scratch.Formula = "=""507.74"""
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Formula = "=""154.46"""
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Formula = "=""2.591.92"""
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Formula = "=""6.56"""
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Formula = "=""3.042.56"""
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Formula = "=""60.525.65"""
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$scratch.Formula = "=""21.62"""
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$scratch.Formula = "=""2.598.08"""
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$scratch.Formula = "=""4.77"""
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Formula = "=""346.10"""
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Formula = "=""10.48"""
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Formula = "=""59.96"""
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Formula = "=""0.420"""
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Formula = "=""0.167"""
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$scratch.Formula = "=""2.702.97"""
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Formula = "=""0.0₃0846"""
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Formula = "=""7.38"""
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Formula = "=""19.38"""
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Formula = "=""152.75"""
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$scratch.Formula = "=""5.71"""
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Formula = "=""3.98"""
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Formula = "=""0.852"""
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Formula = "=""1.48"""
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Formula = "=""0.847"""
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Formula = "=""36.15"""
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Formula = "=""3.75"""
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Formula = "=""296.29"""
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Formula = "=""0.622"""
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Formula = "=""0.0558"""
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Formula = "=""0.997"""
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Formula = "=""19.77"""
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$scratch.Formula = "=""4.88"""
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Clear()

# --- Column E (Volume/1h %) updates ---
# These values keep their padding spaces and percent sign, so Excel already
# treats the assignment as plain text without any extra formatting tricks.
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  -4.60%  "
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("E10").Value = "  +6.54%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  +5.89%  "
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("E43").Value = "  -2.90%  "
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("E50").Value = "  -2.60%  "
